$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F19 4938 -> 4943, F21 829 -> 830
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F19").Value = 4943
$ws1.Range("F21").Value = 830

# Sheet "演出" (sheet2): F2 80 -> 81
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 81

# Sheet "全部类型" (sheet4): F19 4939 -> 4943, F20 80 -> 81, F23 829 -> 830
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F19").Value = 4943
$ws4.Range("F20").Value = 81
$ws4.Range("F23").Value = 830
